# The DEBIT/CREDIT columns (D and E) had been accidentally swapped, and the
# original_index column (F) needs to be dropped from the sheet.
#
# Fix: swap the contents of columns D and E back (header labels included),
# then remove column F entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$rngD = $ws.Range("D1:D$lastRow")
$rngE = $ws.Range("E1:E$lastRow")

# Read both columns first, then write them back swapped, so neither write
# clobbers data the other read still needs.
$valuesD = $rngD.Value()
$valuesE = $rngE.Value()

$rngD.Value = $valuesE
$rngE.Value = $valuesD

# Drop the now-unwanted original_index column.
$ws.Columns("F").Delete()
